# modify row-title judging standard
#
# The worksheet "工作表1" holds a small plot/species table. A new row is
# inserted (pushing the second plot's rows down by one), the plot-title
# cells that used to live in column B (merged with the area header) are
# replaced by explicit row labels "A1"/"A2", the species names are
# changed from the placeholder letters A/B/C to real tree species, and
# the merged header cell that used to only carry formatting now also
# carries the column title text "冠幅".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表1")

# Insert a new blank row above what is currently row 4; this shifts the
# old rows 4-9 down to rows 5-10 (and extends the sheet dimension).
$ws.Rows.Item(4).Insert()

# The old "样地1" plot-title cell (B2) is no longer used - the row title
# now lives in column B of the data rows instead, so clear it.
$ws.Range("B2").ClearContents()

# The merged F2:H2 header cell used to be blank (formatting only); it now
# shows the column title "冠幅" as well.
$ws.Range("F2").Value = "冠幅"

# New row title for the first plot (row 4, freshly inserted).
$ws.Range("B4").Value = "A1"

# Species names for the first plot's rows (now rows 5-7).
$ws.Range("B5").Value = "松树"
$ws.Range("B6").Value = "杨树"
$ws.Range("B7").Value = "柳树"

# Row title for the second plot (previously "样地2" at row 7, now row 8).
$ws.Range("B8").Value = "A2"

# Species names for the second plot's rows (now rows 9-10).
$ws.Range("B9").Value = "松树"
$ws.Range("B10").Value = "柳树"

# Match the author's final selection.
[void]$ws.Range("C7").Select()
